$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: first save of the new staff member ---
$ws.Range("A3").Value = "NATHAN"
$ws.Range("B3").Value = "DANKSIN"
$ws.Range("C3").Value = 2222
$ws.Range("D3").Value = "MANAGER"

# --- Row 4: second save of the same staff member (the "exists" check is
#     missing, so the record gets written again, this time with the
#     corrected surname and the id typed into a text field) ---
$ws.Range("A4").Value = "NATHAN"
$ws.Range("B4").Value = "DANSKIN"

# The staff id is stored as text this time. Writing "2222" straight into
# Value would be auto-recognised as a number, so stage it as text in a
# scratch cell (quote-prefixed so it is kept literal), copy it, and paste
# only the value into C4 - this keeps the cell's default (unstyled) look
# while still storing a text cell.
$ws.Range("F1").Value = "'2222"
$ws.Range("F1").Copy()
$ws.Range("C4").PasteSpecial(-4163)
$ws.Range("F1").Clear()

$ws.Range("D4").Value = "MANAGER"

# --- column width tweaks ---
# (the host rounds ColumnWidth to the nearest 1/6 of a character, so these
# land on the closest representable width to the authored 8.424911 / 10.139196)
$ws.Columns.Item(1).ColumnWidth = 7.666666666666667
$ws.Columns.Item(4).ColumnWidth = 9.333333333333334

# --- selection moves to the newly entered id cell ---
$ws.Range("D4").Select()
